$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value (as text, matching original inlineStr format)
$updates = @{
    2  = @{ D = "61.021.85"; E = "  +7.54%  " }
    3  = @{ D = "3.354.99";  E = "  +3.94%  " }
    4  = @{ D = "1.00";      E = "  -0.15%  " }
    5  = @{ D = "413.66";    E = "  +5.26%  " }
    6  = @{ D = "113.00";    E = "  +6.19%  " }
    7  = @{ E = "  +4.49%  " }
    8  = @{ E = "  -0.10%  " }
    9  = @{ E = "  +4.52%  " }
    10 = @{ D = "39.86";     E = "  +3.36%  " }
    11 = @{ D = "0.1000";    E = "  +4.20%  " }
    12 = @{ D = "0.144";     E = "  +1.44%  " }
    13 = @{ D = "3.887.98";  E = "  +3.59%  " }
    14 = @{ E = "  +4.63%  " }
    15 = @{ D = "19.81";     E = "  +5.11%  " }
    16 = @{ D = "3.357.34";  E = "  +4.16%  " }
    17 = @{ E = "  +2.66%  " }
    18 = @{ D = "60.840.96"; E = "  +7.35%  " }
    19 = @{ D = "10.75";     E = "  +1.71%  " }
    20 = @{ E = "  +3.23%  " }
    21 = @{ D = "0.0000112"; E = "  +6.77%  " }
    22 = @{ D = "12.98";     E = "  +0.97%  " }
    23 = @{ D = "303.64";    E = "  +2.77%  " }
    24 = @{ D = "75.05";     E = "  +2.27%  " }
    25 = @{ E = "  +2.59%  " }
    26 = @{ D = "28.74";     E = "  +3.92%  " }
    27 = @{ D = "4.49";      E = "  +2.66%  " }
    28 = @{ D = "0.180";     E = "  +6.94%  " }
    29 = @{ D = "7.96";      E = "  +3.32%  " }
    30 = @{ D = "7.65";      E = "  +6.67%  " }
    31 = @{ D = "2.63";      E = "  +24.51%  " }
    32 = @{ E = "  +5.49%  " }
    33 = @{ D = "11.45";     E = "  +5.15%  " }
    34 = @{ D = "1.00";      E = "  +0.07%  " }
    35 = @{ D = "39.11";     E = "  +4.76%  " }
    36 = @{ D = "0.0511";    E = "  +6.49%  " }
    37 = @{ D = "52.52";     E = "  +1.90%  " }
    38 = @{ D = "3.12";      E = "  +2.23%  " }
    39 = @{ E = "  -0.11%  " }
    40 = @{ D = "3.41";      E = "  -1.25%  " }
    41 = @{ D = "136.71";    E = "  +2.53%  " }
    42 = @{ D = "0.296";     E = "  +4.23%  " }
    43 = @{ E = "  +3.58%  " }
    44 = @{ D = "1.91";      E = "  +1.48%  " }
    45 = @{ E = "  +1.18%  " }
    46 = @{ D = "16.88";     E = "  -0.16%  " }
    47 = @{
        B = "WEMIXToken"
        C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
        D = "2.25"
        E = "  +8.81%  "
    }
    48 = @{
        B = "EnergySwap"
        C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
        D = "22.51"
        E = "  +3.79%  "
    }
    49 = @{ D = "2.173.74"; E = "  +2.26%  " }
    50 = @{ E = "  +1.36%  " }
    51 = @{ D = "1.97";     E = "  -0.99%  " }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $cellRef = "$colLetter$rowNum"
        $range = $ws.Range($cellRef)
        # Force text interpretation so numeric-looking strings (e.g. "1.00",
        # "113.00", "61.021.85") are preserved exactly as typed, then reset
        # the cell style back to Normal so no stray style gets attached.
        $range.NumberFormat = "@"
        $range.Value = $cols[$colLetter]
        $range.Style = "Normal"
    }
}
